$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row total correct-answer marking value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the "Total" row correct marks total (B12: 45 -> 75)
$ws.Range("B12").Value = 75

# Update the correct/total marks label (E12: "38/84" -> "75/140")
$ws.Range("E12").Value = "75/140"
